$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Active Status" column, header in I1 (mirrors the existing header row style)
$ws.Range("I1").Value = "Active Status"

# Move the selection to the new column, matching the updated layout
$ws.Range("I2").Select() | Out-Null
